$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell H1 with same formatting as G1 (bold, border, centered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Add data cell H2
$ws.Range("H2").Value = 1
